# Apply edits described by the commit:
# "Diagramas de robustez y secuencia CU 10 y 11" — update the status
# ("Estado") and effort ("Esfuerzo") columns for the use cases in rows
# 10, 13, 14 and 15 of the "Casos de Uso" sheet: mark them as
# "planificado" (was "vacio") and set effort to 1 (was 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Rows whose Estado / Esfuerzo values need to change.
$rows = @(10, 13, 14, 15)

foreach ($r in $rows) {
    $ws.Range("E$r").Value = "planificado"
    $ws.Range("F$r").Value = 1
}

# Update the active cell selection on the sheet to match the saved view.
$ws.Activate()
$ws.Range("E14").Select()

$wb.Save()
